$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45178 -> 45179) for every data row (rows 2 through 458).
$lastRow = 458
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value = 45179
    }
}

# Row 4 (case "A 45009-2019") additionally got a new species (Knärot) recorded,
# which bumps several counts and adds a new hyperlink column (U).
$ws.Cells.Item(4, 8).Value = 4    # H4 Fridlysta
$ws.Cells.Item(4, 11).Value = 1   # K4 VU
$ws.Cells.Item(4, 15).Value = 3   # O4 Rödlistade
$ws.Cells.Item(4, 16).Value = 1   # P4 Hotade
$ws.Cells.Item(4, 17).Value = 5   # Q4 Alla arter

$ws.Cells.Item(4, 18).Value = "Knärot`r`nLunglav`r`nTalltita`r`nKorallrot`r`nSpindelblomster"

$ws.Cells.Item(4, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_OVANAKER/knärot/A 45009-2019.png")'
